# feat: implement class managers table (#169)
#
# Updates the attendance-list title year, swaps in a new student
# (HAR JING DARYL) for the first entry of the Class Group A21 roster,
# and corrects a handful of student ID numbers across the A26 / A32
# rosters. Also restores the sheet's scroll position / selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Title: report year 2022 -> 2023
$ws.Range("A1").Value = "Class Attendance List: 2023, 2  Date: 15-JUN-2023 13:01"

# 2. Class Group A21 roster (rows 13-40): student #1 replaced
$ws.Range("B13").Value = "HAR JING DARYL"
$ws.Range("F13").Value = "HARJ0002"

# 3. Class Group A26 roster (rows 47-71): VMS Acc / ID corrections
$ws.Range("F48").Value = "CHUJ6788"
$ws.Range("F51").Value = "GOHK4568"
$ws.Range("F65").Value = "TANK4322"
$ws.Range("F67").Value = "TAYL7655"

# 4. Class Group A32 roster (rows 78-90): VMS Acc / ID corrections
$ws.Range("F80").Value = "GOHK4569"
$ws.Range("F81").Value = "KOHM2346"
$ws.Range("F84").Value = "LIMH5679"

# 5. Restore view state: scroll position + active selection
$win = $excel.ActiveWindow
$win.ScrollRow = 58
$win.ScrollColumn = 1
$ws.Range("F81").Select()
